# Peru Liga 1 - base update (06-04-2024 15:39)
# Applies:
#   1) Odds-row data corrections for matches 175/177 and 186/187/188
#      (the match id in col A and fixed metadata stayed put, only the
#      bookmaker data moved between rows).
#   2) Minor odds corrections on upcoming-fixture rows 273-275.
#   3) A brand new upcoming-fixture row (276) appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 175 (was row 177's data) ---
$ws.Cells.Item(175, 2).Value2 = 7302796
$ws.Cells.Item(175, 6).Value2 = 'Sport Huancayo'
$ws.Cells.Item(175, 7).Value2 = 'Sport Boys'
$ws.Cells.Item(175, 8).Value2 = 1
$ws.Cells.Item(175, 9).Value2 = 0
$ws.Cells.Item(175, 10).Value2 = 'H'
$ws.Cells.Item(175, 11).Value2 = 1.727
$ws.Cells.Item(175, 12).Value2 = 3.75
$ws.Cells.Item(175, 13).Value2 = 4.333
$ws.Cells.Item(175, 14).Value2 = 1.25
$ws.Cells.Item(175, 15).Value2 = 5.25
$ws.Cells.Item(175, 16).Value2 = 10
$ws.Cells.Item(175, 17).Value2 = -1.75
$ws.Cells.Item(175, 18).Value2 = 1.925
$ws.Cells.Item(175, 19).Value2 = 1.875
$ws.Cells.Item(175, 20).Value2 = 3
$ws.Cells.Item(175, 21).Value2 = 1.875
$ws.Cells.Item(175, 22).Value2 = 1.925
$ws.Cells.Item(175, 23).Value2 = 0.25
$ws.Cells.Item(175, 24).Value2 = -1
$ws.Cells.Item(175, 25).Value2 = -1
$ws.Cells.Item(175, 26).Value2 = -1
$ws.Cells.Item(175, 27).Value2 = 0.875
$ws.Cells.Item(175, 28).Value2 = -1
$ws.Cells.Item(175, 29).Value2 = 0.925

# --- Row 177 (was row 175's data) ---
$ws.Cells.Item(177, 2).Value2 = 7302200
$ws.Cells.Item(177, 6).Value2 = 'Carlos Manucci'
$ws.Cells.Item(177, 7).Value2 = 'Deportivo Binacional'
$ws.Cells.Item(177, 8).Value2 = 3
$ws.Cells.Item(177, 9).Value2 = 2
$ws.Cells.Item(177, 10).Value2 = 'H'
$ws.Cells.Item(177, 11).Value2 = 2
$ws.Cells.Item(177, 12).Value2 = 3.2
$ws.Cells.Item(177, 13).Value2 = 3.75
$ws.Cells.Item(177, 14).Value2 = 1.75
$ws.Cells.Item(177, 15).Value2 = 3.4
$ws.Cells.Item(177, 16).Value2 = 4.333
$ws.Cells.Item(177, 17).Value2 = -0.5
$ws.Cells.Item(177, 18).Value2 = 1.85
$ws.Cells.Item(177, 19).Value2 = 1.95
$ws.Cells.Item(177, 20).Value2 = 2.5
$ws.Cells.Item(177, 21).Value2 = 1.85
$ws.Cells.Item(177, 22).Value2 = 1.95
$ws.Cells.Item(177, 23).Value2 = 0.75
$ws.Cells.Item(177, 24).Value2 = -1
$ws.Cells.Item(177, 25).Value2 = -1
$ws.Cells.Item(177, 26).Value2 = 0.8500000000000001
$ws.Cells.Item(177, 27).Value2 = -1
$ws.Cells.Item(177, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(177, 29).Value2 = -1

# --- Row 186 (was row 187's data) ---
$ws.Cells.Item(186, 2).Value2 = 7384629
$ws.Cells.Item(186, 6).Value2 = 'Deportivo Garcilaso'
$ws.Cells.Item(186, 7).Value2 = 'Alianza Lima'
$ws.Cells.Item(186, 8).Value2 = 0
$ws.Cells.Item(186, 9).Value2 = 1
$ws.Cells.Item(186, 10).Value2 = 'A'
$ws.Cells.Item(186, 11).Value2 = 2.625
$ws.Cells.Item(186, 12).Value2 = 3.3
$ws.Cells.Item(186, 13).Value2 = 2.5
$ws.Cells.Item(186, 14).Value2 = 2.7
$ws.Cells.Item(186, 15).Value2 = 3.4
$ws.Cells.Item(186, 16).Value2 = 2.375
$ws.Cells.Item(186, 17).Value2 = 0
$ws.Cells.Item(186, 18).Value2 = 2.025
$ws.Cells.Item(186, 19).Value2 = 1.775
$ws.Cells.Item(186, 20).Value2 = 2.25
$ws.Cells.Item(186, 21).Value2 = 1.825
$ws.Cells.Item(186, 22).Value2 = 1.975
$ws.Cells.Item(186, 23).Value2 = -1
$ws.Cells.Item(186, 24).Value2 = -1
$ws.Cells.Item(186, 25).Value2 = 1.375
$ws.Cells.Item(186, 26).Value2 = -1
$ws.Cells.Item(186, 27).Value2 = 0.7749999999999999
$ws.Cells.Item(186, 28).Value2 = -1
$ws.Cells.Item(186, 29).Value2 = 0.9750000000000001

# --- Row 187 (was row 188's data) ---
$ws.Cells.Item(187, 2).Value2 = 7384625
$ws.Cells.Item(187, 6).Value2 = 'AD Tarma'
$ws.Cells.Item(187, 7).Value2 = 'Carlos Manucci'
$ws.Cells.Item(187, 8).Value2 = 0
$ws.Cells.Item(187, 9).Value2 = 0
$ws.Cells.Item(187, 10).Value2 = 'D'
$ws.Cells.Item(187, 11).Value2 = 1.5
$ws.Cells.Item(187, 12).Value2 = 3.75
$ws.Cells.Item(187, 13).Value2 = 7
$ws.Cells.Item(187, 14).Value2 = 1.363
$ws.Cells.Item(187, 15).Value2 = 4.333
$ws.Cells.Item(187, 16).Value2 = 9.5
$ws.Cells.Item(187, 17).Value2 = -1.25
$ws.Cells.Item(187, 18).Value2 = 1.875
$ws.Cells.Item(187, 19).Value2 = 1.925
$ws.Cells.Item(187, 20).Value2 = 2.5
$ws.Cells.Item(187, 21).Value2 = 1.8
$ws.Cells.Item(187, 22).Value2 = 2
$ws.Cells.Item(187, 23).Value2 = -1
$ws.Cells.Item(187, 24).Value2 = 3.333
$ws.Cells.Item(187, 25).Value2 = -1
$ws.Cells.Item(187, 26).Value2 = -1
$ws.Cells.Item(187, 27).Value2 = 0.925
$ws.Cells.Item(187, 28).Value2 = -1
$ws.Cells.Item(187, 29).Value2 = 1

# --- Row 188 (was row 186's data) ---
$ws.Cells.Item(188, 2).Value2 = 7384630
$ws.Cells.Item(188, 6).Value2 = 'Atletico Grau'
$ws.Cells.Item(188, 7).Value2 = 'Unin Comercio'
$ws.Cells.Item(188, 8).Value2 = 0
$ws.Cells.Item(188, 9).Value2 = 1
$ws.Cells.Item(188, 10).Value2 = 'A'
$ws.Cells.Item(188, 11).Value2 = 2.8
$ws.Cells.Item(188, 12).Value2 = 3.4
$ws.Cells.Item(188, 13).Value2 = 2.15
$ws.Cells.Item(188, 14).Value2 = 1.75
$ws.Cells.Item(188, 15).Value2 = 3.6
$ws.Cells.Item(188, 16).Value2 = 3.8
$ws.Cells.Item(188, 17).Value2 = -0.75
$ws.Cells.Item(188, 18).Value2 = 2
$ws.Cells.Item(188, 19).Value2 = 1.8
$ws.Cells.Item(188, 20).Value2 = 3
$ws.Cells.Item(188, 21).Value2 = 1.85
$ws.Cells.Item(188, 22).Value2 = 1.95
$ws.Cells.Item(188, 23).Value2 = -1
$ws.Cells.Item(188, 24).Value2 = -1
$ws.Cells.Item(188, 25).Value2 = 2.8
$ws.Cells.Item(188, 26).Value2 = -1
$ws.Cells.Item(188, 27).Value2 = 0.8
$ws.Cells.Item(188, 28).Value2 = -1
$ws.Cells.Item(188, 29).Value2 = 0.95

# --- Row 273: small odds refresh on oddAHH / oddAHA ---
$ws.Cells.Item(273, 21).Value2 = 1.95
$ws.Cells.Item(273, 22).Value2 = 1.9

# --- Row 274: small odds refresh ---
$ws.Cells.Item(274, 14).Value2 = 2.6
$ws.Cells.Item(274, 15).Value2 = 3.4
$ws.Cells.Item(274, 16).Value2 = 2.375
$ws.Cells.Item(274, 18).Value2 = 2.05
$ws.Cells.Item(274, 19).Value2 = 1.8

# --- Row 275: small odds refresh ---
$ws.Cells.Item(275, 14).Value2 = 4.5
$ws.Cells.Item(275, 16).Value2 = 1.666
$ws.Cells.Item(275, 17).Value2 = 0.75
$ws.Cells.Item(275, 18).Value2 = 1.9
$ws.Cells.Item(275, 19).Value2 = 1.95
$ws.Cells.Item(275, 21).Value2 = 2.05
$ws.Cells.Item(275, 22).Value2 = 1.8

# --- New row 276: upcoming fixture, appended below row 275 ---
# Copy formatting (bold/border for col A, date format for col E) from the
# row above so the new row matches the sheet's existing look.
$ws.Cells.Item(275, 1).Copy()
$ws.Cells.Item(276, 1).PasteSpecial(-4122)
$ws.Cells.Item(275, 5).Copy()
$ws.Cells.Item(276, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(276, 1).Value2 = 274
$ws.Cells.Item(276, 2).Value2 = 8011511
$ws.Cells.Item(276, 3).Value2 = 'Peru Liga 1'
$ws.Cells.Item(276, 4).Value2 = 'Peru Liga 1'
$ws.Cells.Item(276, 5).Value2 = 45389.72916666666
$ws.Cells.Item(276, 6).Value2 = 'Sporting Cristal'
$ws.Cells.Item(276, 7).Value2 = 'Sport Huancayo'
$ws.Cells.Item(276, 11).Value2 = 1.333
$ws.Cells.Item(276, 12).Value2 = 4.5
$ws.Cells.Item(276, 13).Value2 = 9
$ws.Cells.Item(276, 14).Value2 = 1.4
$ws.Cells.Item(276, 15).Value2 = 4.333
$ws.Cells.Item(276, 16).Value2 = 7
$ws.Cells.Item(276, 17).Value2 = -1.25
$ws.Cells.Item(276, 18).Value2 = 1.9
$ws.Cells.Item(276, 19).Value2 = 1.95
$ws.Cells.Item(276, 20).Value2 = 3
$ws.Cells.Item(276, 21).Value2 = 2.025
$ws.Cells.Item(276, 22).Value2 = 1.825
$ws.Cells.Item(276, 23).Value2 = 0
$ws.Cells.Item(276, 24).Value2 = 0
$ws.Cells.Item(276, 25).Value2 = 0
$ws.Cells.Item(276, 26).Value2 = 0
$ws.Cells.Item(276, 27).Value2 = 0
